# Update "想去人数" (column F) figures across the three sheets that carry
# event data, matching the refreshed output generated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 473
$ws1.Range("F6").Value = 230
$ws1.Range("F7").Value = 212
$ws1.Range("F8").Value = 243
$ws1.Range("F9").Value = 2828
$ws1.Range("F12").Value = 2170
$ws1.Range("F13").Value = 266
$ws1.Range("F19").Value = 1255
$ws1.Range("F20").Value = 4540
$ws1.Range("F22").Value = 4796
$ws1.Range("F23").Value = 1329
$ws1.Range("F24").Value = 2778
$ws1.Range("F25").Value = 3189
$ws1.Range("F26").Value = 138
$ws1.Range("F27").Value = 1485
$ws1.Range("F28").Value = 235
$ws1.Range("F29").Value = 819
$ws1.Range("F30").Value = 86
$ws1.Range("F31").Value = 243
$ws1.Range("F32").Value = 864
$ws1.Range("F33").Value = 1518
$ws1.Range("F35").Value = 238
$ws1.Range("F36").Value = 599
$ws1.Range("F38").Value = 289
$ws1.Range("F39").Value = 360

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 49
$ws2.Range("F7").Value = 47
$ws2.Range("F10").Value = 18

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 473
$ws4.Range("F7").Value = 49
$ws4.Range("F8").Value = 230
$ws4.Range("F9").Value = 212
$ws4.Range("F10").Value = 47
$ws4.Range("F11").Value = 243
$ws4.Range("F12").Value = 2828
$ws4.Range("F15").Value = 2170
$ws4.Range("F16").Value = 266
$ws4.Range("F20").Value = 18
$ws4.Range("F21").Value = 2513
$ws4.Range("F22").Value = 1255
$ws4.Range("F26").Value = 4540
$ws4.Range("F28").Value = 4796
$ws4.Range("F29").Value = 1329
$ws4.Range("F30").Value = 2778
$ws4.Range("F31").Value = 3189
$ws4.Range("F32").Value = 138
$ws4.Range("F35").Value = 1485
$ws4.Range("F37").Value = 235
$ws4.Range("F38").Value = 819
$ws4.Range("F39").Value = 86
$ws4.Range("F40").Value = 243
$ws4.Range("F41").Value = 864
$ws4.Range("F43").Value = 1518
$ws4.Range("F45").Value = 238
$ws4.Range("F46").Value = 599
$ws4.Range("F48").Value = 289
$ws4.Range("F49").Value = 360
